$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add header cells for new columns I0 (I) and IF (J), matching the style of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I0 (column I) and IF (column J), rows 2-55
$data = @(
    @(2, 7, 7),
    @(3, 5, 5),
    @(4, 6, 6),
    @(5, 6, 6),
    @(6, 6, 6),
    @(7, 5, 5),
    @(8, 6, 6),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 6, 6),
    @(12, 9, 9),
    @(13, 7, 8),
    @(14, 8, 8),
    @(15, 6, 6),
    @(16, 8, 8),
    @(17, 6, 6),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 6, 6),
    @(21, 7, 7),
    @(22, 9, 9),
    @(23, 8, 8),
    @(24, 7, 7),
    @(25, 8, 8),
    @(26, 7, 7),
    @(27, 7, 7),
    @(28, 8, 8),
    @(29, 6, 6),
    @(30, 6, 7),
    @(31, 8, 8),
    @(32, 7, 7),
    @(33, 8, 8),
    @(34, 6, 6),
    @(35, 8, 8),
    @(36, 8, 9),
    @(37, 8, 8),
    @(38, 8, 8),
    @(39, 8, 8),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 9, 9),
    @(43, 8, 8),
    @(44, 8, 8),
    @(45, 8, 8),
    @(46, 6, 6),
    @(47, 7, 7),
    @(48, 7, 7),
    @(49, 5, 6),
    @(50, 8, 8),
    @(51, 4, 4),
    @(52, 4, 4),
    @(53, 4, 4),
    @(54, 2, 2),
    @(55, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
